{"js": "// Replace the date line and every division-problem cell text in document\n// order, preserving existing run/paragraph formatting by using\n// search + insertText(\"Replace\") instead of clearing/rewriting the runs.\n\nconst replacements = [\n  [\"2024-06-07 Friday\", \"2024-06-08 Saturday\"],\n  [\"58\u00f75=\", \"22\u00f75=\"],\n  [\"73\u00f78=\", \"83\u00f73=\"],\n  [\"31\u00f77=\", \"11\u00f79=\"],\n  [\"52\u00f73=\", \"64\u00f78=\"],\n  [\"48\u00f79=\", \"59\u00f77=\"],\n  [\"78\u00f78=\", \"95\u00f79=\"],\n  [\"73\u00f78=\", \"57\u00f73=\"],\n  [\"43\u00f73=\", \"10\u00f76=\"],\n  [\"93\u00f78=\", \"44\u00f73=\"],\n  [\"73\u00f74=\", \"40\u00f78=\"],\n  [\"58\u00f72=\", \"45\u00f77=\"],\n  [\"86\u00f73=\", \"33\u00f76=\"],\n  [\"39\u00f72=\", \"49\u00f77=\"],\n  [\"70\u00f72=\", \"21\u00f74=\"],\n  [\"58\u00f79=\", \"93\u00f73=\"],\n  [\"61\u00f73=\", \"50\u00f78=\"],\n  [\"95\u00f74=\", \"33\u00f75=\"],\n  [\"33\u00f72=\", \"92\u00f76=\"],\n  [\"37\u00f73=\", \"59\u00f74=\"],\n  [\"78\u00f79=\", \"42\u00f79=\"],\n  [\"88\u00f77=\", \"23\u00f79=\"],\n  [\"39\u00f72=\", \"43\u00f73=\"],\n  [\"58\u00f72=\", \"16\u00f76=\"],\n  [\"14\u00f73=\", \"74\u00f75=\"],\n  [\"56\u00f73=\", \"56\u00f77=\"],\n];\n\nconst body = context.document.body;\n\n// Issue one search per distinct \"from\" string (search results come back\n// in document order), then consume the Nth occurrence the Nth time that\n// \"from\" string shows up in our ordered replacement list.\nconst uniqueFroms = [...new Set(replacements.map((pair) => pair[0]))];\nconst searchResults = {};\nfor (const from of uniqueFroms) {\n  searchResults[from] = body.search(from, { matchCase: true });\n  searchResults[from].load(\"items\");\n}\nawait context.sync();\n\nconst nextOccurrenceIndex = {};\nfor (const [from, to] of replacements) {\n  const idx = nextOccurrenceIndex[from] || 0;\n  nextOccurrenceIndex[from] = idx + 1;\n  const range = searchResults[from].items[idx];\n  range.insertText(to, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Replace the date line and every division-problem cell text, in document\n# order, preserving existing run/paragraph formatting.\n#\n# We reuse a single Range/Find pair anchored at the very start of the\n# document. Each wdReplaceOne execution both performs the substitution and\n# advances the range's search position past the replaced text, so issuing\n# the calls in document order naturally lands on the correct occurrence\n# even when the same source text (e.g. \"73\u00f78=\") appears more than once.\n\n$d = $word.ActiveDocument\n$rng = $d.Content\n$rng.Start = 0\n$find = $rng.Find\n\n$pairs = @(\n    @(\"2024-06-07 Friday\", \"2024-06-08 Saturday\"),\n    @(\"58\u00f75=\", \"22\u00f75=\"),\n    @(\"73\u00f78=\", \"83\u00f73=\"),\n    @(\"31\u00f77=\", \"11\u00f79=\"),\n    @(\"52\u00f73=\", \"64\u00f78=\"),\n    @(\"48\u00f79=\", \"59\u00f77=\"),\n    @(\"78\u00f78=\", \"95\u00f79=\"),\n    @(\"73\u00f78=\", \"57\u00f73=\"),\n    @(\"43\u00f73=\", \"10\u00f76=\"),\n    @(\"93\u00f78=\", \"44\u00f73=\"),\n    @(\"73\u00f74=\", \"40\u00f78=\"),\n    @(\"58\u00f72=\", \"45\u00f77=\"),\n    @(\"86\u00f73=\", \"33\u00f76=\"),\n    @(\"39\u00f72=\", \"49\u00f77=\"),\n    @(\"70\u00f72=\", \"21\u00f74=\"),\n    @(\"58\u00f79=\", \"93\u00f73=\"),\n    @(\"61\u00f73=\", \"50\u00f78=\"),\n    @(\"95\u00f74=\", \"33\u00f75=\"),\n    @(\"33\u00f72=\", \"92\u00f76=\"),\n    @(\"37\u00f73=\", \"59\u00f74=\"),\n    @(\"78\u00f79=\", \"42\u00f79=\"),\n    @(\"88\u00f77=\", \"23\u00f79=\"),\n    @(\"39\u00f72=\", \"43\u00f73=\"),\n    @(\"58\u00f72=\", \"16\u00f76=\"),\n    @(\"14\u00f73=\", \"74\u00f75=\"),\n    @(\"56\u00f73=\", \"56\u00f77=\")\n)\n\nforeach ($pair in $pairs) {\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 0, $false, [ref]$find.Replacement.Text, 1) | Out-Null\n}\n"}
